$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B64").Value = "Film de fusor negro HP"
$ws.Range("D64").Value = 15000
$ws.Range("E64").Value = 100000
$ws.Range("F64").Value = 49
$ws.Range("G64").Value = 64
$ws.Range("J64").Value = 735000
